# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
#
# The "Estado de Cuenta" detail table (rows 16-31) is rebuilt: instead of
# being grouped by period with workers interleaved, the rows are now
# grouped by worker (JORGE ARTURO MARTINEZ VASQUEZ first, then ELIANA
# MENDEZ BLANCO), each with periods 2208 down to 2201, and ELIANA's
# "Salario Basico" (column G) is updated from 908526 to 1423500.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$tipoDoc = "CC"

$jorgeDoc    = "73135739"
$jorgeNombre = "JORGE ARTURO MARTINEZ VASQUEZ"
$jorgeSalario = 908526

$elianaDoc    = "1143324046"
$elianaNombre = "ELIANA MENDEZ BLANCO"
$elianaSalario = 1423500

# row -> (doc, nombre, periodo, valorMora, salarioBasico)
$rows = @(
    @{ R = 16; Doc = $jorgeDoc;  Nombre = $jorgeNombre;  Periodo = "2208"; Mora = 25749; Salario = $jorgeSalario },
    @{ R = 17; Doc = $jorgeDoc;  Nombre = $jorgeNombre;  Periodo = "2207"; Mora = 35112; Salario = $jorgeSalario },
    @{ R = 18; Doc = $jorgeDoc;  Nombre = $jorgeNombre;  Periodo = "2206"; Mora = 35112; Salario = $jorgeSalario },
    @{ R = 19; Doc = $jorgeDoc;  Nombre = $jorgeNombre;  Periodo = "2205"; Mora = 35112; Salario = $jorgeSalario },
    @{ R = 20; Doc = $jorgeDoc;  Nombre = $jorgeNombre;  Periodo = "2204"; Mora = 35112; Salario = $jorgeSalario },
    @{ R = 21; Doc = $jorgeDoc;  Nombre = $jorgeNombre;  Periodo = "2203"; Mora = 35112; Salario = $jorgeSalario },
    @{ R = 22; Doc = $jorgeDoc;  Nombre = $jorgeNombre;  Periodo = "2202"; Mora = 35112; Salario = $jorgeSalario },
    @{ R = 23; Doc = $jorgeDoc;  Nombre = $jorgeNombre;  Periodo = "2201"; Mora = 35112; Salario = $jorgeSalario },
    @{ R = 24; Doc = $elianaDoc; Nombre = $elianaNombre; Periodo = "2208"; Mora = 25749; Salario = $elianaSalario },
    @{ R = 25; Doc = $elianaDoc; Nombre = $elianaNombre; Periodo = "2207"; Mora = 36341; Salario = $elianaSalario },
    @{ R = 26; Doc = $elianaDoc; Nombre = $elianaNombre; Periodo = "2206"; Mora = 36341; Salario = $elianaSalario },
    @{ R = 27; Doc = $elianaDoc; Nombre = $elianaNombre; Periodo = "2205"; Mora = 36341; Salario = $elianaSalario },
    @{ R = 28; Doc = $elianaDoc; Nombre = $elianaNombre; Periodo = "2204"; Mora = 36341; Salario = $elianaSalario },
    @{ R = 29; Doc = $elianaDoc; Nombre = $elianaNombre; Periodo = "2203"; Mora = 36341; Salario = $elianaSalario },
    @{ R = 30; Doc = $elianaDoc; Nombre = $elianaNombre; Periodo = "2202"; Mora = 36341; Salario = $elianaSalario },
    @{ R = 31; Doc = $elianaDoc; Nombre = $elianaNombre; Periodo = "2201"; Mora = 36341; Salario = $elianaSalario }
)

foreach ($row in $rows) {
    $r = $row.R
    $ws.Range("B$r").Value = $tipoDoc
    $ws.Range("C$r").Value = $row.Doc
    $ws.Range("D$r").Value = $row.Nombre
    $ws.Range("E$r").Value = $row.Periodo
    $ws.Range("F$r").Value = $row.Mora
    $ws.Range("G$r").Value = $row.Salario
}
